$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2022.03.31"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = 4.92
$ws.Range("C8").Value = 4.213
$ws.Range("D8").Value = 24014
$ws.Range("E8").Value = 5700

$ws.Range("B9").Select()
